$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 177 (shifts existing rows 177-191 down to 178-192)
$ws.Rows.Item(177).Insert()

# Fill in the new row 177 with data
$ws.Cells.Item(177, 1).Value = 5
$ws.Cells.Item(177, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(177, 3).Value = "Maule"
$ws.Cells.Item(177, 4).Value = 44461
$ws.Cells.Item(177, 5).Value = 7
$ws.Cells.Item(177, 6).Value = 100112023
$ws.Cells.Item(177, 7).Value = "Brócoli"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 3000
$ws.Cells.Item(177, 11).Value = 600
$ws.Cells.Item(177, 12).Value = 600
$ws.Cells.Item(177, 13).Value = 600
$ws.Cells.Item(177, 14).Value = "`$/unidad"
$ws.Cells.Item(177, 15).Value = "Región del Maule"
$ws.Cells.Item(177, 16).Value = 600
$ws.Cells.Item(177, 17).Value = 1
$ws.Cells.Item(177, 18).Value = "Hortaliza"

# Apply the date style (same as other D cells) to the new D177 cell
$ws.Cells.Item(176, 4).Copy()
$ws.Cells.Item(177, 4).PasteSpecial(-4122)
$ws.Cells.Item(177, 4).Value = 44461
